# Update the "Förändrad" (Changed) date column (C) for data rows 2-28
# from serial date 45551 (2024-09-16) to 45552 (2024-09-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDateSerial = 45552

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDateSerial
}
